$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 724, shifting existing rows 724+ down by one
# (matches the diff: rows 724-765 become 725-766, dimension grows to D766).
$ws.Rows.Item(724).Insert()

# Populate the newly inserted row with the new data point.
# Force column A to stay a plain text value (it looks like a date, so the
# auto-detect would otherwise coerce it into a date serial + date format).
$ws.Cells.Item(724, 1).NumberFormat = "@"
$ws.Cells.Item(724, 1).Value = "2026/01/30"
$ws.Cells.Item(724, 1).Style = "Normal"

$ws.Cells.Item(724, 2).Value = "金"
$ws.Cells.Item(724, 3).Value = 17
$ws.Cells.Item(724, 4).Value = 20
